# Insert two new weekly price records (rows) right before the current row 33
# (Cilantro, Terminal Hortofrutícola Agro Chillán), pushing all subsequent
# rows down by two. The new rows share the same Mercado / Región / Categoría /
# Variedad / Unidad / Origen / Kg-o-Unidades / Clasificación metadata as their
# neighbours, only the date, quality, volume and price columns differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 33 (each Insert() pushes existing row 33
# downward, so doing it twice opens up rows 33 and 34).
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# --- New row 33: Primera, 13-07-2022 ---------------------------------------
$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(33, 3).Value = "Ñuble"
$ws.Cells.Item(33, 4).Value = 44755
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 5).Value = 16
$ws.Cells.Item(33, 6).Value = 100112040
$ws.Cells.Item(33, 7).Value = "Cilantro"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 200
$ws.Cells.Item(33, 11).Value = 600
$ws.Cells.Item(33, 12).Value = 700
$ws.Cells.Item(33, 13).Value = 650
$ws.Cells.Item(33, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(33, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(33, 16).Value = 650
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# --- New row 34: Segunda, 13-07-2022 ---------------------------------------
$ws.Cells.Item(34, 1).Value = 7
$ws.Cells.Item(34, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(34, 3).Value = "Ñuble"
$ws.Cells.Item(34, 4).Value = 44755
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 16
$ws.Cells.Item(34, 6).Value = 100112040
$ws.Cells.Item(34, 7).Value = "Cilantro"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 150
$ws.Cells.Item(34, 11).Value = 500
$ws.Cells.Item(34, 12).Value = 500
$ws.Cells.Item(34, 13).Value = 500
$ws.Cells.Item(34, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(34, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(34, 16).Value = 500
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
